$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 163, pushing existing rows 163:271 down to 164:272
$ws.Rows.Item(163).Insert()

# Populate the newly inserted row 163 with the new weekly record
$ws.Cells.Item(163, 1).Value = 7
$ws.Cells.Item(163, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(163, 3).Value = 'Ñuble'
$ws.Cells.Item(163, 4).Value = (Get-Date -Year 2022 -Month 9 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(163, 5).Value = 16
$ws.Cells.Item(163, 6).Value = 100112003
$ws.Cells.Item(163, 7).Value = 'Ajo'
$ws.Cells.Item(163, 8).Value = 'Chino'
$ws.Cells.Item(163, 9).Value = 'Primera'
$ws.Cells.Item(163, 10).Value = 40
$ws.Cells.Item(163, 11).Value = 22000
$ws.Cells.Item(163, 12).Value = 23000
$ws.Cells.Item(163, 13).Value = 22500
$ws.Cells.Item(163, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(163, 15).Value = 'China'
$ws.Cells.Item(163, 16).Value = 2250
$ws.Cells.Item(163, 17).Value = 10
$ws.Cells.Item(163, 18).Value = 'Hortaliza'
